# Generate Report for Handoff
# The localization status moved from "In Translation" to "Ready for handoff"
# and the Xliff generation / handoff timestamps were refreshed. Updating the
# cell text makes Excel recompute the (auto-fitted) width of the "Status"
# columns to accommodate the longer "Ready for handoff" string.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status + latest generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 01:04:28"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 01:04:23"

# --- de-de sheet: Status (Latest Handoff Datetime shares the Overview date) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 01:04:28"

# Widen the Status columns to fit the new, longer "Ready for handoff" text
# (mirrors Excel's auto-fit-on-edit behaviour for these report columns).
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
